# Updates cryptos list: price/volume refresh plus a few row re-orderings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is plain text (coin names / URLs) or a numeric-looking
# string that Excel will NOT silently reinterpret as a number (because it has
# more than one '.' or other non-numeric characters). These can be set
# directly via .Value without any extra bookkeeping.
$plainMap = @{
    "D2"  = "47.787.46"
    "E2"  = "  +1.41%  "
    "D3"  = "2.499.49"
    "E3"  = "  +0.54%  "
    "E5"  = "  -0.17%  "
    "E6"  = "  +2.05%  "
    "E9"  = "  +1.99%  "
    "E10" = "  +5.88%  "
    "E11" = "  -0.18%  "
    "B12" = "Chainlink"
    "C12" = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
    "E12" = "  +3.00%  "
    "B13" = "TRON"
    "C13" = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
    "E13" = "  +0.47%  "
    "E14" = "  +0.24%  "
    "D15" = "2.889.88"
    "E15" = "  +0.45%  "
    "D16" = "2.496.97"
    "E16" = "  +0.90%  "
    "E17" = "  -0.53%  "
    "D18" = "47.677.14"
    "E18" = "  +1.30%  "
    "E19" = "  +2.09%  "
    "E20" = "  -0.29%  "
    "E21" = "  +12.60%  "
    "D22" = "0.0₃0942"
    "E22" = "  +0.29%  "
    "E23" = "  +0.15%  "
    "E24" = "  -0.93%  "
    "E25" = "  -1.31%  "
    "E26" = "  +0.20%  "
    "E27" = "  -1.07%  "
    "E28" = "  -0.49%  "
    "B29" = "Kaspa"
    "C29" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "E29" = "  +1.23%  "
    "B30" = "InjectiveProtocol"
    "C30" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "E30" = "  +0.34%  "
    "B31" = "Toncoin"
    "C31" = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
    "E31" = "  -9.46%  "
    "E32" = "  +0.92%  "
    "E33" = "  +0.71%  "
    "E34" = "  -2.24%  "
    "E35" = "  -0.68%  "
    "E36" = "  +0.04%  "
    "E37" = "  -1.09%  "
    "E38" = "  -1.01%  "
    "E39" = "  -0.75%  "
    "E40" = "  -0.30%  "
    "E41" = "  +4.98%  "
    "E42" = "  -1.39%  "
    "E43" = "  -2.43%  "
    "E44" = "  -0.30%  "
    "D45" = "2.004.90"
    "E45" = "  +1.81%  "
    "E46" = "  +1.78%  "
    "E48" = "  +0.20%  "
    "E49" = "  -0.42%  "
    "E50" = "  -2.56%  "
    "E51" = "  +2.70%  "
}

foreach ($key in $plainMap.Keys) {
    $ws.Range($key).Value = $plainMap[$key]
}

# Cells whose new value is a "clean" numeric-looking string (single decimal
# point, digits only). Left alone, Excel's COM layer would coerce these into
# a true number cell (changing both the cell type and, via NumberFormat,
# the style index). To keep them as text - matching the original inline
# string cells exactly - force a text number format before assigning the
# value, then restore the cell to the default "Normal" style so no stray
# style index is left behind.
$textNumMap = @{
    "D5"  = "322.80"
    "D6"  = "108.97"
    "D9"  = "0.551"
    "D10" = "40.28"
    "D12" = "19.00"
    "D13" = "0.124"
    "D17" = "0.849"
    "D19" = "13.13"
    "D23" = "70.77"
    "D24" = "247.94"
    "D29" = "0.139"
    "D30" = "35.00"
    "D31" = "2.08"
    "D32" = "49.91"
    "D33" = "19.85"
    "D34" = "5.35"
    "D35" = "0.0788"
    "D36" = "1.00"
    "D41" = "22.27"
    "D43" = "119.47"
    "D49" = "9.02"
}

foreach ($key in $textNumMap.Keys) {
    $rng = $ws.Range($key)
    $rng.NumberFormat = "@"
    $rng.Value = $textNumMap[$key]
    $rng.Style = "Normal"
}
